# Apply KHL injuries snapshot refresh:
#  - snapshot sheet: remove 3 players who returned from injury
#    (АМР Гиздатуллин Артур; СЕВ Ващенко Григорий; СЕВ Грудинин Владимир)
#    and refresh the scraped_at timestamp for every remaining row.
#  - returned sheet: append the 3 players that just returned.

$wb = $excel.ActiveWorkbook
$snapshot = $wb.Worksheets.Item("snapshot")
$returned = $wb.Worksheets.Item("returned")

# --- 1. Remove returned players from the snapshot sheet -------------------
# Row numbers (1-based, as currently laid out in the sheet) of the players
# that are leaving the injury list. Delete bottom-up so earlier indices stay
# valid while we work.
$snapshot.Rows.Item(22).Delete() | Out-Null   # СЕВ Грудинин Владимир
$snapshot.Rows.Item(21).Delete() | Out-Null   # СЕВ Ващенко Григорий
$snapshot.Rows.Item(10).Delete() | Out-Null   # АМР Гиздатуллин Артур

# --- 2. Refresh scraped_at (column K) for every remaining data row --------
$newScrapedAt = @(
    "2025-11-13T09:46:51.454806+00:00",
    "2025-11-13T09:46:51.454829+00:00",
    "2025-11-13T09:46:51.454887+00:00",
    "2025-11-13T09:46:54.045003+00:00",
    "2025-11-13T09:46:54.045042+00:00",
    "2025-11-13T09:46:54.045065+00:00",
    "2025-11-13T09:46:56.296170+00:00",
    "2025-11-13T09:46:58.570210+00:00",
    "2025-11-13T09:47:01.414786+00:00",
    "2025-11-13T09:47:01.414857+00:00",
    "2025-11-13T09:47:06.077465+00:00",
    "2025-11-13T09:47:08.385171+00:00",
    "2025-11-13T09:47:11.266096+00:00",
    "2025-11-13T09:47:11.266130+00:00",
    "2025-11-13T09:47:11.266151+00:00",
    "2025-11-13T09:47:13.626142+00:00",
    "2025-11-13T09:47:16.482744+00:00",
    "2025-11-13T09:47:16.482779+00:00",
    "2025-11-13T09:47:19.333270+00:00",
    "2025-11-13T09:47:21.674157+00:00",
    "2025-11-13T09:47:21.674188+00:00",
    "2025-11-13T09:47:21.674207+00:00",
    "2025-11-13T09:47:21.674224+00:00",
    "2025-11-13T09:47:21.674241+00:00",
    "2025-11-13T09:47:24.023834+00:00",
    "2025-11-13T09:47:24.023865+00:00",
    "2025-11-13T09:47:26.264440+00:00",
    "2025-11-13T09:47:26.264475+00:00",
    "2025-11-13T09:47:26.264495+00:00",
    "2025-11-13T09:47:29.031770+00:00",
    "2025-11-13T09:47:29.031798+00:00",
    "2025-11-13T09:47:29.031815+00:00",
    "2025-11-13T09:47:31.243083+00:00",
    "2025-11-13T09:47:31.243112+00:00",
    "2025-11-13T09:47:31.243129+00:00",
    "2025-11-13T09:47:31.243144+00:00",
    "2025-11-13T09:47:31.243163+00:00",
    "2025-11-13T09:47:31.243179+00:00",
    "2025-11-13T09:47:34.060908+00:00",
    "2025-11-13T09:47:34.060939+00:00",
    "2025-11-13T09:47:38.774881+00:00",
    "2025-11-13T09:47:38.774918+00:00",
    "2025-11-13T09:47:38.774940+00:00",
    "2025-11-13T09:47:38.774959+00:00",
    "2025-11-13T09:47:41.488374+00:00",
    "2025-11-13T09:47:41.488405+00:00"
)

for ($i = 0; $i -lt $newScrapedAt.Length; $i++) {
    $row = $i + 2
    $snapshot.Cells.Item($row, 11).Value = $newScrapedAt[$i]
}

# --- 3. Append the returned players to the "returned" sheet ---------------
# Note: changed_day ("2025-11-13") is prefixed with a leading apostrophe so
# Excel stores it as plain text instead of auto-converting it to a date
# serial number (matches how the existing changed_day cells are stored).
$returnedRows = @(
    @("АМР", "Амур", "Гиздатуллин Артур", "1369_АМР_гиздатуллинартур", "RETURN", "2025-11-13T17:47:41.991047+08:00", "'2025-11-13"),
    @("СЕВ", "Северсталь", "Ващенко Григорий", "1369_СЕВ_ващенкогригорий", "RETURN", "2025-11-13T17:47:41.991047+08:00", "'2025-11-13"),
    @("СЕВ", "Северсталь", "Грудинин Владимир", "1369_СЕВ_грудининвладимир", "RETURN", "2025-11-13T17:47:41.991047+08:00", "'2025-11-13")
)

$startRow = 4
for ($i = 0; $i -lt $returnedRows.Length; $i++) {
    $r = $startRow + $i
    $data = $returnedRows[$i]
    for ($c = 0; $c -lt $data.Length; $c++) {
        $returned.Cells.Item($r, $c + 1).Value = $data[$c]
    }
}
